# Update the "取得日時" (retrieved datetime) timestamps in the "ランサーズ" sheet
# for all existing data rows (2-10) from 2026-01-19 12:44:23 to 2026-01-19 12:58:24.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2026-01-19 12:44:23"
$newTimestamp = "2026-01-19 12:58:24"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Text -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
